$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 8).Value = 1455265.2
$ws.Cells.Item(2, 10).Value = 2002
$ws.Cells.Item(2, 12).Value = 2002
$ws.Cells.Item(2, 14).Value = -2228
$ws.Cells.Item(12, 8).Value = 325
$ws.Cells.Item(12, 9).Value = 325
$ws.Cells.Item(12, 10).Value = 0
$ws.Cells.Item(12, 11).Value = 325
$ws.Cells.Item(12, 12).Value = 0
$ws.Cells.Item(12, 13).Value = -155
$ws.Cells.Item(12, 14).ClearContents()
$ws.Cells.Item(28, 8).Value = 831.2727
$ws.Cells.Item(28, 9).Value = 494.1579
$ws.Cells.Item(28, 11).Value = 494.1579
$ws.Cells.Item(28, 13).Value = -9.157899999999984
$ws.Cells.Item(32, 8).Value = 2343.8333
$ws.Cells.Item(32, 9).Value = 2000
$ws.Cells.Item(32, 10).Value = 2375.0908
$ws.Cells.Item(32, 11).Value = 2000
$ws.Cells.Item(32, 12).Value = 2375.0908
$ws.Cells.Item(32, 14).Value = -3027.0908
$ws.Cells.Item(32, 13).Value = -1674
$ws.Cells.Item(38, 8).Value = 942.2857
$ws.Cells.Item(38, 9).Value = 116.7
$ws.Cells.Item(38, 10).Value = 3006.25
$ws.Cells.Item(38, 11).Value = 350.1
$ws.Cells.Item(38, 12).Value = 9018.75
$ws.Cells.Item(38, 13).Value = 21.89999999999998
$ws.Cells.Item(38, 14).Value = -9762.75
$ws.Cells.Item(39, 8).Value = 1466.6957
$ws.Cells.Item(39, 9).Value = 854
$ws.Cells.Item(39, 10).Value = 3202.6667
$ws.Cells.Item(39, 11).Value = 2562
$ws.Cells.Item(39, 12).Value = 9608.000100000001
$ws.Cells.Item(39, 13).Value = -2266
$ws.Cells.Item(39, 14).Value = -10200.0001
$ws.Cells.Item(43, 8).Value = 378.625
$ws.Cells.Item(43, 9).Value = 344.5
$ws.Cells.Item(43, 11).Value = 344.5
$ws.Cells.Item(43, 13).Value = -275.5
$ws.Cells.Item(98, 8).Value = 8029.231
$ws.Cells.Item(98, 9).Value = 1198.3334
$ws.Cells.Item(98, 11).Value = 1198.3334
$ws.Cells.Item(98, 13).Value = 299.6666
$ws.Cells.Item(103, 8).Value = 1500.3125
$ws.Cells.Item(103, 9).Value = 1695
$ws.Cells.Item(103, 10).Value = 1250
$ws.Cells.Item(103, 11).Value = 5085
$ws.Cells.Item(103, 12).Value = 3750
$ws.Cells.Item(103, 13).Value = -4499
$ws.Cells.Item(103, 14).Value = -4922
$ws.Cells.Item(122, 8).Value = 8029.231
$ws.Cells.Item(122, 9).Value = 1198.3334
$ws.Cells.Item(122, 11).Value = 3595.0002
$ws.Cells.Item(122, 13).Value = -1145.0002
$ws.Cells.Item(125, 8).Value = 7000
$ws.Cells.Item(125, 10).Value = 7000
$ws.Cells.Item(125, 12).Value = 63000
$ws.Cells.Item(125, 14).Value = -67920
$ws.Cells.Item(128, 8).Value = 149998
$ws.Cells.Item(128, 10).Value = 149998
$ws.Cells.Item(128, 12).Value = 149998
$ws.Cells.Item(128, 14).Value = -159958
$ws.Cells.Item(132, 8).Value = 3285.3518
$ws.Cells.Item(132, 9).Value = 2910.6326
$ws.Cells.Item(132, 10).Value = 6957.6
$ws.Cells.Item(132, 11).Value = 8731.897799999999
$ws.Cells.Item(132, 12).Value = 20872.8
$ws.Cells.Item(132, 13).Value = -6201.897799999999
$ws.Cells.Item(132, 14).Value = -25932.8
$ws.Cells.Item(133, 8).Value = 90639.28999999999
$ws.Cells.Item(133, 10).Value = 90639.28999999999
$ws.Cells.Item(133, 12).Value = 90639.28999999999
$ws.Cells.Item(133, 14).Value = -100759.29
$ws.Cells.Item(136, 8).Value = 96259.336
$ws.Cells.Item(136, 10).Value = 96259.336
$ws.Cells.Item(136, 12).Value = 96259.336
$ws.Cells.Item(136, 14).Value = -106459.336
$ws.Cells.Item(137, 8).Value = 1171193.5
$ws.Cells.Item(137, 9).Value = 1378.5
$ws.Cells.Item(137, 10).Value = 7410207
$ws.Cells.Item(137, 11).Value = 4135.5
$ws.Cells.Item(137, 12).Value = 22230621
$ws.Cells.Item(137, 13).Value = -1585.5
$ws.Cells.Item(137, 14).Value = -22235721
$ws.Cells.Item(138, 8).Value = 3136
$ws.Cells.Item(138, 9).Value = 2133.1765
$ws.Cells.Item(138, 10).Value = 3910.9092
$ws.Cells.Item(138, 11).Value = 6399.529500000001
$ws.Cells.Item(138, 12).Value = 11732.7276
$ws.Cells.Item(138, 13).Value = -1259.529500000001
$ws.Cells.Item(138, 14).Value = -22012.7276

# Sheet: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 8).Value = 987.5714
$ws.Cells.Item(2, 9).Value = 1005.62964
$ws.Cells.Item(2, 10).Value = 500
$ws.Cells.Item(2, 11).Value = 1005.62964
$ws.Cells.Item(2, 12).Value = 500
$ws.Cells.Item(2, 13).Value = -892.62964
$ws.Cells.Item(2, 14).Value = -726
$ws.Cells.Item(19, 8).Value = 0
$ws.Cells.Item(19, 9).Value = 0
$ws.Cells.Item(19, 11).Value = 0
$ws.Cells.Item(19, 13).ClearContents()
$ws.Cells.Item(25, 8).Value = 11672
$ws.Cells.Item(25, 10).Value = 20000
$ws.Cells.Item(25, 12).Value = 20000
$ws.Cells.Item(25, 14).Value = -20804
$ws.Cells.Item(27, 8).Value = 0
$ws.Cells.Item(27, 10).Value = 0
$ws.Cells.Item(27, 12).Value = 0
$ws.Cells.Item(27, 14).ClearContents()
$ws.Cells.Item(32, 8).Value = 16926.328
$ws.Cells.Item(32, 9).Value = 19480.223
$ws.Cells.Item(32, 10).Value = 3135.3
$ws.Cells.Item(32, 11).Value = 19480.223
$ws.Cells.Item(32, 12).Value = 3135.3
$ws.Cells.Item(32, 13).Value = -19193.223
$ws.Cells.Item(32, 14).Value = -3709.3
$ws.Cells.Item(61, 8).Value = 2089.125
$ws.Cells.Item(61, 9).Value = 1847.7858
$ws.Cells.Item(61, 11).Value = 1847.7858
$ws.Cells.Item(61, 13).Value = -1635.7858
$ws.Cells.Item(63, 8).Value = 2846.5454
$ws.Cells.Item(63, 9).Value = 2145.889
$ws.Cells.Item(63, 11).Value = 2145.889
$ws.Cells.Item(63, 13).Value = -1459.889
$ws.Cells.Item(66, 8).Value = 2846.5454
$ws.Cells.Item(66, 9).Value = 2145.889
$ws.Cells.Item(66, 11).Value = 10729.445
$ws.Cells.Item(66, 13).Value = -7297.445
$ws.Cells.Item(74, 8).Value = 1372.7755
$ws.Cells.Item(74, 9).Value = 809.675
$ws.Cells.Item(74, 10).Value = 3875.4443
$ws.Cells.Item(74, 11).Value = 809.675
$ws.Cells.Item(74, 12).Value = 3875.4443
$ws.Cells.Item(74, 13).Value = 64.32500000000005
$ws.Cells.Item(74, 14).Value = -5623.4443
$ws.Cells.Item(77, 8).Value = 1372.7755
$ws.Cells.Item(77, 9).Value = 809.675
$ws.Cells.Item(77, 10).Value = 3875.4443
$ws.Cells.Item(77, 11).Value = 4048.375
$ws.Cells.Item(77, 12).Value = 19377.2215
$ws.Cells.Item(77, 13).Value = 319.625
$ws.Cells.Item(77, 14).Value = -28113.2215
$ws.Cells.Item(97, 8).Value = 7017.2354
$ws.Cells.Item(97, 10).Value = 1573.3334
$ws.Cells.Item(97, 12).Value = 1573.3334
$ws.Cells.Item(97, 14).Value = -2565.3334
$ws.Cells.Item(102, 8).Value = 3659.5454
$ws.Cells.Item(102, 9).Value = 3700.5
$ws.Cells.Item(102, 11).Value = 3700.5
$ws.Cells.Item(102, 13).Value = -2078.5
$ws.Cells.Item(116, 8).Value = 987.5714
$ws.Cells.Item(116, 9).Value = 1005.62964
$ws.Cells.Item(116, 10).Value = 500
$ws.Cells.Item(116, 11).Value = 1005.62964
$ws.Cells.Item(116, 12).Value = 500
$ws.Cells.Item(116, 13).Value = 1288.37036
$ws.Cells.Item(116, 14).Value = -5088
$ws.Cells.Item(122, 8).Value = 1600.2142
$ws.Cells.Item(122, 9).Value = 1327.0476
$ws.Cells.Item(122, 11).Value = 3981.142800000001
$ws.Cells.Item(122, 13).Value = -1531.142800000001
$ws.Cells.Item(132, 8).Value = 26379.441
$ws.Cells.Item(132, 9).Value = 32498.576
$ws.Cells.Item(132, 11).Value = 97495.728
$ws.Cells.Item(132, 13).Value = -94965.728
$ws.Cells.Item(136, 8).Value = 2089.125
$ws.Cells.Item(136, 9).Value = 1847.7858
$ws.Cells.Item(136, 11).Value = 5543.357400000001
$ws.Cells.Item(136, 13).Value = -2993.357400000001

# Sheet: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(3, 8).Value = 987.5714
$ws.Cells.Item(3, 9).Value = 1005.62964
$ws.Cells.Item(3, 10).Value = 500
$ws.Cells.Item(3, 11).Value = 1005.62964
$ws.Cells.Item(3, 12).Value = 500
$ws.Cells.Item(3, 13).Value = -891.62964
$ws.Cells.Item(3, 14).Value = -728
$ws.Cells.Item(11, 8).Value = 957
$ws.Cells.Item(11, 9).Value = 559.8
$ws.Cells.Item(11, 11).Value = 559.8
$ws.Cells.Item(11, 13).Value = -419.8
$ws.Cells.Item(22, 8).Value = 125561
$ws.Cells.Item(22, 9).Value = 125561
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 11).Value = 125561
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 13).Value = -125388
$ws.Cells.Item(22, 14).ClearContents()
$ws.Cells.Item(68, 8).Value = 250104220
$ws.Cells.Item(68, 10).Value = 333427970
$ws.Cells.Item(68, 12).Value = 333427970
$ws.Cells.Item(68, 14).Value = -333429592
$ws.Cells.Item(71, 8).Value = 250104220
$ws.Cells.Item(71, 10).Value = 333427970
$ws.Cells.Item(71, 12).Value = 1000283910
$ws.Cells.Item(71, 14).Value = -1000292022
$ws.Cells.Item(94, 8).Value = 4107.1
$ws.Cells.Item(94, 9).Value = 3532
$ws.Cells.Item(94, 11).Value = 3532
$ws.Cells.Item(94, 13).Value = -3081
$ws.Cells.Item(99, 8).Value = 73014.8
$ws.Cells.Item(99, 9).Value = 129104.25
$ws.Cells.Item(99, 11).Value = 129104.25
$ws.Cells.Item(99, 13).Value = -127606.25
$ws.Cells.Item(105, 8).Value = 4379.0713
$ws.Cells.Item(105, 9).Value = 4150.9
$ws.Cells.Item(105, 11).Value = 4150.9
$ws.Cells.Item(105, 13).Value = -2403.9
$ws.Cells.Item(134, 8).Value = 1636.025
$ws.Cells.Item(134, 9).Value = 1633.2703
$ws.Cells.Item(134, 11).Value = 4899.810899999999
$ws.Cells.Item(134, 13).Value = -2364.810899999999

# Sheet: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(7, 8).Value = 361.66666
$ws.Cells.Item(7, 9).Value = 482.5
$ws.Cells.Item(7, 11).Value = 482.5
$ws.Cells.Item(7, 13).Value = -369.5
$ws.Cells.Item(31, 8).Value = 1638.9445
$ws.Cells.Item(31, 9).Value = 1693.2941
$ws.Cells.Item(31, 10).Value = 715
$ws.Cells.Item(31, 11).Value = 1693.2941
$ws.Cells.Item(31, 12).Value = 715
$ws.Cells.Item(31, 13).Value = -1398.2941
$ws.Cells.Item(31, 14).Value = -1305
$ws.Cells.Item(34, 8).Value = 1638.9445
$ws.Cells.Item(34, 9).Value = 1693.2941
$ws.Cells.Item(34, 10).Value = 715
$ws.Cells.Item(34, 11).Value = 1693.2941
$ws.Cells.Item(34, 12).Value = 715
$ws.Cells.Item(34, 13).Value = -1491.2941
$ws.Cells.Item(34, 14).Value = -1119
$ws.Cells.Item(58, 8).Value = 78968.16
$ws.Cells.Item(58, 9).Value = 92553.37
$ws.Cells.Item(58, 11).Value = 92553.37
$ws.Cells.Item(58, 13).Value = -92350.37
$ws.Cells.Item(94, 8).Value = 1192.238
$ws.Cells.Item(94, 10).Value = 1288.3846
$ws.Cells.Item(94, 12).Value = 1288.3846
$ws.Cells.Item(94, 14).Value = -2190.3846
$ws.Cells.Item(132, 8).Value = 1727.1666
$ws.Cells.Item(132, 9).Value = 1548.3077
$ws.Cells.Item(132, 10).Value = 2889.75
$ws.Cells.Item(132, 11).Value = 4644.9231
$ws.Cells.Item(132, 12).Value = 8669.25
$ws.Cells.Item(132, 13).Value = -2114.9231
$ws.Cells.Item(132, 14).Value = -13729.25
$ws.Cells.Item(134, 8).Value = 45405.78
$ws.Cells.Item(134, 9).Value = 51826.6
$ws.Cells.Item(134, 10).Value = 2600.3333
$ws.Cells.Item(134, 11).Value = 155479.8
$ws.Cells.Item(134, 12).Value = 7800.999899999999
$ws.Cells.Item(134, 13).Value = -152944.8
$ws.Cells.Item(134, 14).Value = -12870.9999
$ws.Cells.Item(136, 8).Value = 78968.16
$ws.Cells.Item(136, 9).Value = 92553.37
$ws.Cells.Item(136, 11).Value = 277660.11
$ws.Cells.Item(136, 13).Value = -275110.11

# Sheet: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(29, 8).Value = 109.333336
$ws.Cells.Item(29, 9).Value = 114
$ws.Cells.Item(29, 10).Value = 100
$ws.Cells.Item(29, 11).Value = 342
$ws.Cells.Item(29, 12).Value = 300
$ws.Cells.Item(29, 13).Value = -65
$ws.Cells.Item(29, 14).Value = -854
$ws.Cells.Item(44, 8).Value = 2624.6667
$ws.Cells.Item(44, 9).Value = 99.2
$ws.Cells.Item(44, 10).Value = 4428.5713
$ws.Cells.Item(44, 11).Value = 297.6
$ws.Cells.Item(44, 12).Value = 13285.7139
$ws.Cells.Item(44, 13).Value = 100.4
$ws.Cells.Item(44, 14).Value = -14081.7139
$ws.Cells.Item(46, 8).Value = 737.5454999999999
$ws.Cells.Item(46, 9).Value = 752.4706
$ws.Cells.Item(46, 11).Value = 2257.4118
$ws.Cells.Item(46, 13).Value = -2166.4118
$ws.Cells.Item(107, 8).Value = 385599.53
$ws.Cells.Item(107, 10).Value = 401003.53
$ws.Cells.Item(107, 12).Value = 1203010.59
$ws.Cells.Item(107, 14).Value = -1206850.59
$ws.Cells.Item(131, 8).Value = 9766.1875
$ws.Cells.Item(131, 9).Value = 9951
$ws.Cells.Item(131, 11).Value = 29853
$ws.Cells.Item(131, 13).Value = -24813
$ws.Cells.Item(139, 8).Value = 6300.4443
$ws.Cells.Item(139, 9).Value = 902
$ws.Cells.Item(139, 11).Value = 2706
$ws.Cells.Item(139, 13).Value = 2434
$ws.Cells.Item(140, 8).Value = 3122.7222
$ws.Cells.Item(140, 9).Value = 2883.9375
$ws.Cells.Item(140, 11).Value = 8651.8125
$ws.Cells.Item(140, 13).Value = -3471.8125

# Sheet: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(7, 8).Value = 3500500
$ws.Cells.Item(7, 9).Value = 1500
$ws.Cells.Item(7, 10).Value = 5250000
$ws.Cells.Item(7, 11).Value = 1500
$ws.Cells.Item(7, 12).Value = 5250000
$ws.Cells.Item(7, 13).Value = -1388
$ws.Cells.Item(7, 14).Value = -5250224
$ws.Cells.Item(8, 8).Value = 3500500
$ws.Cells.Item(8, 9).Value = 1500
$ws.Cells.Item(8, 10).Value = 5250000
$ws.Cells.Item(8, 11).Value = 1500
$ws.Cells.Item(8, 12).Value = 5250000
$ws.Cells.Item(8, 13).Value = -1361
$ws.Cells.Item(8, 14).Value = -5250278
$ws.Cells.Item(23, 8).Value = 20004.5
$ws.Cells.Item(23, 10).Value = 20004.5
$ws.Cells.Item(23, 12).Value = 20004.5
$ws.Cells.Item(23, 14).Value = -20450.5
$ws.Cells.Item(33, 8).Value = 13299.667
$ws.Cells.Item(33, 10).Value = 14999.5
$ws.Cells.Item(33, 12).Value = 14999.5
$ws.Cells.Item(33, 14).Value = -15503.5
$ws.Cells.Item(80, 8).Value = 7169.9443
$ws.Cells.Item(80, 9).Value = 5880.8
$ws.Cells.Item(80, 10).Value = 8781.375
$ws.Cells.Item(80, 11).Value = 5880.8
$ws.Cells.Item(80, 12).Value = 8781.375
$ws.Cells.Item(80, 13).Value = -4882.8
$ws.Cells.Item(80, 14).Value = -10777.375
$ws.Cells.Item(83, 8).Value = 7169.9443
$ws.Cells.Item(83, 9).Value = 5880.8
$ws.Cells.Item(83, 10).Value = 8781.375
$ws.Cells.Item(83, 11).Value = 29404
$ws.Cells.Item(83, 12).Value = 43906.875
$ws.Cells.Item(83, 13).Value = -24412
$ws.Cells.Item(83, 14).Value = -53890.875
$ws.Cells.Item(97, 8).Value = 3028.625
$ws.Cells.Item(97, 9).Value = 2700.5
$ws.Cells.Item(97, 10).Value = 3356.75
$ws.Cells.Item(97, 11).Value = 2700.5
$ws.Cells.Item(97, 12).Value = 3356.75
$ws.Cells.Item(97, 13).Value = -2204.5
$ws.Cells.Item(97, 14).Value = -4348.75
$ws.Cells.Item(102, 8).Value = 5650.7
$ws.Cells.Item(102, 9).Value = 5851.3335
$ws.Cells.Item(102, 10).Value = 5349.75
$ws.Cells.Item(102, 11).Value = 5851.3335
$ws.Cells.Item(102, 12).Value = 5349.75
$ws.Cells.Item(102, 13).Value = -4229.3335
$ws.Cells.Item(102, 14).Value = -8593.75
$ws.Cells.Item(122, 8).Value = 2554.9688
$ws.Cells.Item(122, 9).Value = 1495.0834
$ws.Cells.Item(122, 11).Value = 4485.2502
$ws.Cells.Item(122, 13).Value = -2035.2502
$ws.Cells.Item(126, 8).Value = 5365.846
$ws.Cells.Item(126, 9).Value = 5385.6
$ws.Cells.Item(126, 11).Value = 16156.8
$ws.Cells.Item(126, 13).Value = -13686.8
$ws.Cells.Item(132, 8).Value = 32188.53
$ws.Cells.Item(132, 9).Value = 36149.7
$ws.Cells.Item(132, 11).Value = 108449.1
$ws.Cells.Item(132, 13).Value = -105919.1

# Sheet: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(7, 8).Value = 4292.4707
$ws.Cells.Item(7, 9).Value = 3399
$ws.Cells.Item(7, 10).Value = 4779.8184
$ws.Cells.Item(7, 11).Value = 3399
$ws.Cells.Item(7, 12).Value = 4779.8184
$ws.Cells.Item(7, 13).Value = -3287
$ws.Cells.Item(7, 14).Value = -5003.8184
$ws.Cells.Item(40, 8).Value = 9733.223
$ws.Cells.Item(40, 10).Value = 4732.6665
$ws.Cells.Item(40, 12).Value = 4732.6665
$ws.Cells.Item(40, 14).Value = -5004.6665
$ws.Cells.Item(55, 8).Value = 506.92593
$ws.Cells.Item(55, 9).Value = 284.44446
$ws.Cells.Item(55, 10).Value = 951.8889
$ws.Cells.Item(55, 11).Value = 284.44446
$ws.Cells.Item(55, 12).Value = 951.8889
$ws.Cells.Item(55, 13).Value = -111.44446
$ws.Cells.Item(55, 14).Value = -1297.8889
$ws.Cells.Item(61, 8).Value = 4980.5835
$ws.Cells.Item(61, 9).Value = 4837.1797
$ws.Cells.Item(61, 10).Value = 5602
$ws.Cells.Item(61, 11).Value = 4837.1797
$ws.Cells.Item(61, 12).Value = 5602
$ws.Cells.Item(61, 13).Value = -4635.1797
$ws.Cells.Item(61, 14).Value = -6006
$ws.Cells.Item(68, 8).Value = 5299.846
$ws.Cells.Item(68, 9).Value = 4334
$ws.Cells.Item(68, 10).Value = 5589.6
$ws.Cells.Item(68, 11).Value = 4334
$ws.Cells.Item(68, 12).Value = 5589.6
$ws.Cells.Item(68, 13).Value = -3585
$ws.Cells.Item(68, 14).Value = -7087.6
$ws.Cells.Item(71, 8).Value = 5299.846
$ws.Cells.Item(71, 9).Value = 4334
$ws.Cells.Item(71, 10).Value = 5589.6
$ws.Cells.Item(71, 11).Value = 21670
$ws.Cells.Item(71, 12).Value = 27948
$ws.Cells.Item(71, 13).Value = -17926
$ws.Cells.Item(71, 14).Value = -35436
$ws.Cells.Item(82, 8).Value = 2994.4167
$ws.Cells.Item(82, 9).Value = 861.3333
$ws.Cells.Item(82, 10).Value = 3705.4443
$ws.Cells.Item(82, 11).Value = 861.3333
$ws.Cells.Item(82, 12).Value = 3705.4443
$ws.Cells.Item(82, 13).Value = -500.3333
$ws.Cells.Item(82, 14).Value = -4427.4443
$ws.Cells.Item(85, 8).Value = 2994.4167
$ws.Cells.Item(85, 9).Value = 861.3333
$ws.Cells.Item(85, 10).Value = 3705.4443
$ws.Cells.Item(85, 11).Value = 861.3333
$ws.Cells.Item(85, 12).Value = 3705.4443
$ws.Cells.Item(85, 13).Value = 386.6667
$ws.Cells.Item(85, 14).Value = -6201.4443
$ws.Cells.Item(100, 8).Value = 6479.6875
$ws.Cells.Item(100, 9).Value = 5485.1113
$ws.Cells.Item(100, 10).Value = 7758.4287
$ws.Cells.Item(100, 11).Value = 5485.1113
$ws.Cells.Item(100, 12).Value = 7758.4287
$ws.Cells.Item(100, 13).Value = -4944.1113
$ws.Cells.Item(100, 14).Value = -8840.4287
$ws.Cells.Item(113, 8).Value = 4980.5835
$ws.Cells.Item(113, 9).Value = 4837.1797
$ws.Cells.Item(113, 10).Value = 5602
$ws.Cells.Item(113, 11).Value = 4837.1797
$ws.Cells.Item(113, 12).Value = 5602
$ws.Cells.Item(113, 13).Value = -2667.1797
$ws.Cells.Item(113, 14).Value = -9942
$ws.Cells.Item(122, 8).Value = 6635.0347
$ws.Cells.Item(122, 9).Value = 5409.45
$ws.Cells.Item(122, 10).Value = 9358.556
$ws.Cells.Item(122, 11).Value = 16228.35
$ws.Cells.Item(122, 12).Value = 28075.668
$ws.Cells.Item(122, 13).Value = -13778.35
$ws.Cells.Item(122, 14).Value = -32975.66800000001
$ws.Cells.Item(126, 8).Value = 4292.4707
$ws.Cells.Item(126, 9).Value = 3399
$ws.Cells.Item(126, 10).Value = 4779.8184
$ws.Cells.Item(126, 11).Value = 10197
$ws.Cells.Item(126, 12).Value = 14339.4552
$ws.Cells.Item(126, 13).Value = -7727
$ws.Cells.Item(126, 14).Value = -19279.4552
$ws.Cells.Item(128, 8).Value = 0
$ws.Cells.Item(128, 10).Value = 0
$ws.Cells.Item(128, 12).Value = 0
$ws.Cells.Item(128, 14).ClearContents()
$ws.Cells.Item(132, 8).Value = 73469.7
$ws.Cells.Item(132, 9).Value = 101970.414
$ws.Cells.Item(132, 11).Value = 305911.242
$ws.Cells.Item(132, 13).Value = -303381.242
$ws.Cells.Item(136, 8).Value = 3616.0952
$ws.Cells.Item(136, 9).Value = 2832.3572
$ws.Cells.Item(136, 10).Value = 5183.5713
$ws.Cells.Item(136, 11).Value = 8497.071599999999
$ws.Cells.Item(136, 12).Value = 15550.7139
$ws.Cells.Item(136, 13).Value = -5947.071599999999
$ws.Cells.Item(136, 14).Value = -20650.7139

# Sheet: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(3, 8).Value = 17666
$ws.Cells.Item(3, 9).Value = 999
$ws.Cells.Item(3, 11).Value = 999
$ws.Cells.Item(3, 13).Value = -885
$ws.Cells.Item(11, 8).Value = 29999.5
$ws.Cells.Item(11, 9).Value = 20000
$ws.Cells.Item(11, 11).Value = 20000
$ws.Cells.Item(11, 13).Value = -19858
$ws.Cells.Item(49, 8).Value = 36999
$ws.Cells.Item(49, 9).Value = 29999
$ws.Cells.Item(49, 11).Value = 29999
$ws.Cells.Item(49, 13).Value = -29769
$ws.Cells.Item(95, 8).Value = 25333.334
$ws.Cells.Item(95, 10).Value = 25333.334
$ws.Cells.Item(95, 12).Value = 25333.334
$ws.Cells.Item(95, 14).Value = -30825.334
$ws.Cells.Item(100, 8).Value = 730.86664
$ws.Cells.Item(100, 9).Value = 588.75
$ws.Cells.Item(100, 10).Value = 1299.3334
$ws.Cells.Item(100, 11).Value = 1177.5
$ws.Cells.Item(100, 12).Value = 2598.6668
$ws.Cells.Item(100, 13).Value = -636.5
$ws.Cells.Item(100, 14).Value = -3680.6668
$ws.Cells.Item(113, 8).Value = 3910.4211
$ws.Cells.Item(113, 9).Value = 2954.9092
$ws.Cells.Item(113, 10).Value = 5224.25
$ws.Cells.Item(113, 11).Value = 8864.7276
$ws.Cells.Item(113, 12).Value = 15672.75
$ws.Cells.Item(113, 13).Value = -6694.7276
$ws.Cells.Item(113, 14).Value = -20012.75
$ws.Cells.Item(122, 8).Value = 2792.862
$ws.Cells.Item(122, 9).Value = 2856.5
$ws.Cells.Item(122, 11).Value = 8569.5
$ws.Cells.Item(122, 13).Value = -6119.5
$ws.Cells.Item(126, 8).Value = 4325.452
$ws.Cells.Item(126, 9).Value = 4071.9412
$ws.Cells.Item(126, 11).Value = 12215.8236
$ws.Cells.Item(126, 13).Value = -9745.8236
$ws.Cells.Item(132, 8).Value = 21389.844
$ws.Cells.Item(132, 9).Value = 22089.428
$ws.Cells.Item(132, 11).Value = 66268.284
$ws.Cells.Item(132, 13).Value = -63738.284
$ws.Cells.Item(136, 8).Value = 3076.6191
$ws.Cells.Item(136, 9).Value = 1995.0625
$ws.Cells.Item(136, 11).Value = 5985.1875
$ws.Cells.Item(136, 13).Value = -3435.1875
